$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.870.13"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").Value = "2.674.89"
$ws.Range("E3").Value = "  +7.66%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'113.58"
$ws.Range("E5").Value = "  +8.93%  "
$ws.Range("D6").Value = "'325.94"
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.553"
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "'40.73"
$ws.Range("E10").Value = "  +5.29%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  +3.13%  "
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("E14").Value = "  +4.65%  "
$ws.Range("D15").Value = "3.095.45"
$ws.Range("E15").Value = "  +7.64%  "
$ws.Range("D16").Value = "2.678.33"
$ws.Range("E16").Value = "  +7.08%  "
$ws.Range("E17").Value = "  +6.10%  "
$ws.Range("D18").Value = "49.875.28"
$ws.Range("E18").Value = "  +4.30%  "
$ws.Range("D19").Value = "'13.14"
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("D20").Value = "'6.78"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("D21").Value = "'2.93"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'71.77"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'275.65"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").Value = "'26.84"
$ws.Range("E26").Value = "  +5.06%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = "  +6.63%  "
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").Value = "'36.21"
$ws.Range("E30").Value = "  +5.47%  "
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").Value = "'50.25"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("E33").Value = "  +4.72%  "
$ws.Range("D34").Value = "'19.53"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("D35").Value = "'0.0806"
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  +12.68%  "
$ws.Range("E38").Value = "  +7.46%  "
$ws.Range("D39").Value = "'3.14"
$ws.Range("E39").Value = "  +9.94%  "
$ws.Range("D40").Value = "'125.55"
$ws.Range("E40").Value = "  +4.95%  "
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("D42").Value = "'22.55"
$ws.Range("E42").Value = "  +5.42%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "'0.0318"
$ws.Range("E44").Value = "  +6.74%  "
$ws.Range("D45").Value = "2.123.21"
$ws.Range("E45").Value = "  +6.95%  "
$ws.Range("D46").Value = "'3.32"
$ws.Range("E46").Value = "  +7.11%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.08"
$ws.Range("E47").Value = "  +9.29%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.22"
$ws.Range("E48").Value = "  +7.56%  "
$ws.Range("D49").Value = "'9.04"
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("D51").Value = "'59.25"
$ws.Range("E51").Value = "  +7.18%  "

Write-Host "Applied cryptos update"
